$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New response row appended to the survey data (row 21).
# Column B ("27") is an age value stored as text in the source data, so
# force the Text number format before assigning it to avoid Excel's
# automatic numeric type inference turning it into a Number.
$ws.Range("B21").NumberFormat = "@"

$ws.Range("A21").Value = "2026-02-12 15:44:47"
$ws.Range("B21").Value = "27"
$ws.Range("C21").Value = "male"
$ws.Range("D21").Value = "Sometimes"
$ws.Range("E21").Value = "Never"
$ws.Range("F21").Value = "Once"
$ws.Range("G21").Value = "I know exactly where to go"
$ws.Range("H21").Value = "Somewhat uncomfortable"
$ws.Range("I21").Value = "Disagree"
$ws.Range("J21").Value = "Sometimes"
$ws.Range("K21").Value = "Rarely"
$ws.Range("L21").Value = "Somewhat"
$ws.Range("M21").Value = "Rarely"
$ws.Range("N21").Value = "Neutral"
$ws.Range("O21").Value = "Rarely"
$ws.Range("P21").Value = "I thought about it"
$ws.Range("Q21").Value = "Disagree"
$ws.Range("R21").Value = "A few times"
$ws.Range("S21").Value = "A little"
$ws.Range("T21").Value = "Sometimes"
$ws.Range("U21").Value = "Very uncomfortable"
$ws.Range("V21").Value = "Disagree"
$ws.Range("W21").Value = "Slightly confident"
